$d = $word.ActiveDocument

# Locate the paragraph containing "Gabriel Santana Goes      RA:819220395"
# (the last paragraph in the body before the section break) and append a
# brand-new paragraph right after it, inheriting its run/paragraph
# formatting (Arial, 24 half-points), containing the text "Testes ".

$targetPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$targetRange = $targetPara.Range
$targetRange.Collapse(0)
$targetRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.InsertBefore("Testes ")
